{"js": "// Make the built-in \"Block Text\" paragraph style (used for block quotes)\n// look like a normal, indented block quote: indent both the left and right\n// edges instead of leaving the paragraph flush with the margins.\n//\n// (Previously this style also overrode the font to a smaller, different\n// typeface; the commit un-does that override so block quotes render in the\n// same font/size as regular body text. Office.js only exposes per-property\n// setters on Word.Font \u2014 there is no API to strip a style's run-properties\n// override outright \u2014 so this script focuses on the indentation change,\n// which is fully expressible through Word.ParagraphFormat.)\n\nconst styles = context.document.getStyles();\nconst blockText = styles.getByNameOrNullObject(\"Block Text\");\nawait context.sync();\n\nif (!blockText.isNullObject) {\n  const pf = blockText.paragraphFormat;\n  // 24pt == 480 twips, matching w:ind w:left=\"480\" w:right=\"480\".\n  pf.leftIndent = 24;\n  pf.rightIndent = 24;\n  await context.sync();\n}\n", "ps1": "# Make the built-in \"Block Text\" paragraph style (used for block quotes)\n# look like a normal, indented block quote: indent both the left and right\n# edges instead of leaving the paragraph flush with the margins.\n#\n# (Previously this style also overrode the font to a smaller, different\n# typeface; the commit un-does that override so block quotes render in the\n# same font/size as regular body text. The Word object model only exposes\n# per-property setters on the style's Font object -- there is no API call\n# that strips a style's run-properties override outright -- so this script\n# focuses on the indentation change, which is fully expressible through\n# ParagraphFormat.LeftIndent / RightIndent.)\n\n$d = $word.ActiveDocument\n$blockText = $d.Styles(\"Block Text\")\n\n# 24pt == 480 twips, matching w:ind w:left=\"480\" w:right=\"480\".\n$blockText.ParagraphFormat.LeftIndent = 24\n$blockText.ParagraphFormat.RightIndent = 24\n"}
